$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 296, shifting existing rows 296-370 down to 297-371.
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new data record.
$ws.Cells.Item(296, 1).Value2 = 10
$ws.Cells.Item(296, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(296, 3).Value2 = "La Araucanía"
$ws.Cells.Item(296, 4).Value2 = 44736
$ws.Cells.Item(296, 5).Value2 = 9
$ws.Cells.Item(296, 6).Value2 = 100112037
$ws.Cells.Item(296, 7).Value2 = "Cebollín"
$ws.Cells.Item(296, 8).Value2 = "Sin especificar"
$ws.Cells.Item(296, 9).Value2 = "Primera"
$ws.Cells.Item(296, 10).Value2 = 30
$ws.Cells.Item(296, 11).Value2 = 9000
$ws.Cells.Item(296, 12).Value2 = 9000
$ws.Cells.Item(296, 13).Value2 = 9000
$ws.Cells.Item(296, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(296, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(296, 16).Value2 = 750
$ws.Cells.Item(296, 17).Value2 = 12
$ws.Cells.Item(296, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by the other rows in column D.
$ws.Cells.Item(296, 4).NumberFormat = $ws.Cells.Item(297, 4).NumberFormat
